$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Add the new "Std" / "Relative std" summary columns for the three
#    per-frequency blocks (low / medium / high).
# ---------------------------------------------------------------------------

# Low frequency block (data B2:B31, average lives in E4)
$ws.Range("D15").Value = "Std"
$ws.Range("D15").Font.Bold = $true
$ws.Range("E15").Value = "Relative std"
$ws.Range("E15").Font.Bold = $true
$ws.Range("D16").Formula = "=STDEV(B2:B31)"
$ws.Range("E16").Formula = "= (D16 / E4) * 100"

# Medium frequency block (data B34:B63, average lives in E36)
$ws.Range("D47").Value = "Std"
$ws.Range("D47").Font.Bold = $true
$ws.Range("E47").Value = "Relative std"
$ws.Range("E47").Font.Bold = $true
$ws.Range("D48").Formula = "=STDEV(B34:B63)"
$ws.Range("E48").Formula = "= (D48 / E36) * 100"

# High frequency block (data B66:B95, average lives in E68)
$ws.Range("D79").Value = "Std"
$ws.Range("D79").Font.Bold = $true
$ws.Range("E79").Value = "Relative std"
$ws.Range("E79").Font.Bold = $true
$ws.Range("D80").Formula = "=STDEV(B66:B95)"
$ws.Range("E80").Formula = "= (D80 / E68) * 100"

# ---------------------------------------------------------------------------
# 2) Re-create the hidden "_xlchart.v1.*" defined names in the new order that
#    Excel produced after refreshing/editing the box-whisker charts. This
#    collapses the stray duplicate entries (v1.10 / v1.11) and renumbers the
#    remaining ones.
# ---------------------------------------------------------------------------

$existing = @()
for ($i = 1; $i -le $wb.Names.Count; $i++) {
    $existing += , $wb.Names.Item($i).Name
}
foreach ($nm in $existing) {
    if ($nm -like "_xlchart.v1.*") {
        $wb.Names.Item($nm).Delete()
    }
}

$wb.Names.Add("_xlchart.v1.0", "=Blad1!`$A`$66:`$A`$95")
$wb.Names.Add("_xlchart.v1.1", "=Blad1!`$B`$65")
$wb.Names.Add("_xlchart.v1.2", "=Blad1!`$B`$66:`$B`$95")
$wb.Names.Add("_xlchart.v1.3", "=Blad1!`$A`$98:`$A`$187")
$wb.Names.Add("_xlchart.v1.4", "=Blad1!`$B`$98:`$B`$187")
$wb.Names.Add("_xlchart.v1.5", "=Blad1!`$A`$34:`$A`$63")
$wb.Names.Add("_xlchart.v1.6", "=Blad1!`$B`$33")
$wb.Names.Add("_xlchart.v1.7", "=Blad1!`$B`$34:`$B`$63")
$wb.Names.Add("_xlchart.v1.8", "=Blad1!`$A`$2:`$A`$31")
$wb.Names.Add("_xlchart.v1.9", "=Blad1!`$B`$2:`$B`$31")

# These are the internal "chart data" defined names Excel hides from the
# Name Box / Name Manager UI.
for ($i = 0; $i -le 9; $i++) {
    $wb.Names.Item("_xlchart.v1.$i").Visible = $false
}

# ---------------------------------------------------------------------------
# 3) Update the active selection on the sheet to match where the edit
#    happened (D79:E80, the new high-frequency Std/Relative std cells).
# ---------------------------------------------------------------------------

$ws.Range("D79:E80").Select() | Out-Null
